# Adds two new paragraphs after the "Kaisen : Mejoramiento Continuo." paragraph:
#   1) "KPI : Key Process Indicator" (with spell-check proofErr markers around
#      "Process" and "Indicator", matching Word's auto spell-check behaviour)
#   2) a new empty paragraph with the same run/paragraph formatting
#
# Both new paragraphs reuse the same pPr/rPr (color 000000/themeColor text1,
# lang es-MX, spacing after=100 line=240 auto) already used by the
# surrounding "Concepts" paragraphs in this document.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "Kaisen : Mejoramiento Continuo."
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Mejoramiento Continuo\.") {
        $anchor = $p
    }
}

$rng = $anchor.Range
$rng.Collapse(0)

# Create the two new paragraph marks right after the anchor paragraph; Word
# copies the anchor's pPr/rPr onto them automatically.
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$anchorIdx = $anchor.Index
$p1 = $d.Paragraphs.Item($anchorIdx + 1)
$p2 = $d.Paragraphs.Item($anchorIdx + 2)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="es-MX"/></w:rPr>'
$pPr = '<w:pPr><w:spacing w:after="100" w:line="240" w:lineRule="auto"/>' + $rPr + '</w:pPr>'

# Paragraph 1: "KPI : Key Process Indicator" with spellcheck proofErr
# wrappers around "Process" and "Indicator" (mirrors how Word marks the
# other foreign/unrecognised words - Hoshin, Kanri, Kaisen - elsewhere in
# this document).
$p1xml = '<w:p ' + $wns + '>' + $pPr +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">KPI : Key </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>Process</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>Indicator</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

# Paragraph 2: empty paragraph, same formatting, no runs.
$p2xml = '<w:p ' + $wns + '>' + $pPr + '</w:p>'

[void]$p1.Range.InsertXML($p1xml)
[void]$p2.Range.InsertXML($p2xml)
